$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja2")
$ws.Range("B3:B41").Value = "LE_06_07_CO"
